$wb = $excel.ActiveWorkbook

# ---- Proximity sheet: append rows 41-49 ----
$ws = $wb.Worksheets.Item("Proximity")

$proximityRows = @(
    @("2026-02-01", "15:15:14", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "15:15:16", "15:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "15:15:20", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "15:15:22", "15:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "15:15:27", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "15:15:47", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "15:15:49", "15:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "15:15:59", "15:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "15:16:02", "15:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door")
)

$startRow = 41
for ($i = 0; $i -lt $proximityRows.Count; $i++) {
    $r = $startRow + $i
    $row = $proximityRows[$i]

    # Force column A to be stored as plain text ("2026-02-01" would
    # otherwise be auto-recognized as a date by Excel), then restore the
    # default "Normal" style so no extra formatting is applied.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# ---- Camera sheet: append rows 12-16 ----
$ws2 = $wb.Worksheets.Item("Camera")

$cameraRows = @(
    @("2026-02-01", "15:15:16", "15:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "15:15:22", "15:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "15:15:29", "15:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "15:15:49", "15:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "15:16:02", "15:00", "Living Room Main Door", "Image Captured", "Active")
)

$startRow2 = 12
for ($i = 0; $i -lt $cameraRows.Count; $i++) {
    $r = $startRow2 + $i
    $row = $cameraRows[$i]

    $ws2.Cells.Item($r, 1).NumberFormat = "@"
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 1).Style = "Normal"

    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
}
